$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "culture_collection" column (AM) entirely. This shifts all
# subsequent data/columns left by one, but comment objects stay anchored
# to their original cell addresses, so we must manually re-home the
# comment text for every shifted header cell afterwards.
[void]$ws.Range("AM:AM").EntireColumn.Delete()

# Re-point each header comment (AM15..DA15) to the text that used to
# belong to the next column over, mirroring the column shift.
[void]$ws.Range('AM15').Comment.Text('density of sample')
[void]$ws.Range('AN15').Comment.Text('concentration of diether lipids; can include multiple types of diether lipids')
[void]$ws.Range('AO15').Comment.Text('concentration of dissolved carbon dioxide')
[void]$ws.Range('AP15').Comment.Text('concentration of dissolved hydrogen')
[void]$ws.Range('AQ15').Comment.Text('dissolved inorganic carbon concentration')
[void]$ws.Range('AR15').Comment.Text('concentration of dissolved organic carbon')
[void]$ws.Range('AS15').Comment.Text('dissolved organic nitrogen concentration measured as; total dissolved nitrogen - NH4 - NO3 - NO2')
[void]$ws.Range('AT15').Comment.Text('concentration of dissolved oxygen')
[void]$ws.Range('AU15').Comment.Text('Plasmids that have significance phenotypic consequence')
[void]$ws.Range('AV15').Comment.Text('measurement of glucosidase activity')
[void]$ws.Range('AW15').Comment.Text('Health or disease status of sample at time of collection')
[void]$ws.Range('AX15').Comment.Text('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".')
[void]$ws.Range('AY15').Comment.Text('NCBI taxonomy ID of the host, e.g. 9606')
[void]$ws.Range('AZ15').Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
[void]$ws.Range('BA15').Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
[void]$ws.Range('BB15').Comment.Text('concentration of magnesium')
[void]$ws.Range('BC15').Comment.Text('measurement of mean friction velocity')
[void]$ws.Range('BD15').Comment.Text('measurement of mean peak friction velocity')
[void]$ws.Range('BE15').Comment.Text('methane (gas) amount or concentration at the time of sampling')
[void]$ws.Range('BF15').Comment.Text('any other measurement performed or parameter collected, that is not listed here')
[void]$ws.Range('BG15').Comment.Text('concentration of n-alkanes; can include multiple n-alkanes')
[void]$ws.Range('BH15').Comment.Text('concentration of nitrate')
[void]$ws.Range('BI15').Comment.Text('concentration of nitrite')
[void]$ws.Range('BJ15').Comment.Text('concentration of nitrogen (total)')
[void]$ws.Range('BK15').Comment.Text('concentration of organic carbon')
[void]$ws.Range('BL15').Comment.Text('concentration of organic matter')
[void]$ws.Range('BM15').Comment.Text('concentration of organic nitrogen')
[void]$ws.Range('BN15').Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
[void]$ws.Range('BO15').Comment.Text('oxygenation status of sample')
[void]$ws.Range('BP15').Comment.Text('concentration of particulate organic carbon')
[void]$ws.Range('BQ15').Comment.Text('particles are classified, based on their size, into six general categories: clay, silt, sand, gravel, cobbles, and boulders; should include amount of particle preceded by the name of the particle type; can include multiple values')
[void]$ws.Range('BR15').Comment.Text('To what is the entity pathogenic')
[void]$ws.Range('BS15').Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
[void]$ws.Range('BT15').Comment.Text('concentration of petroleum hydrocarbon')
[void]$ws.Range('BU15').Comment.Text('pH measurement')
[void]$ws.Range('BV15').Comment.Text('concentration of phaeopigments; can include multiple phaeopigments')
[void]$ws.Range('BW15').Comment.Text('concentration of phosphate')
[void]$ws.Range('BX15').Comment.Text('concentration of phospholipid fatty acids; can include multiple values')
[void]$ws.Range('BY15').Comment.Text('porosity of deposited sediment is volume of voids divided by the total volume of sample')
[void]$ws.Range('BZ15').Comment.Text('concentration of potassium')
[void]$ws.Range('CA15').Comment.Text('pressure to which the sample is subject, in atmospheres')
[void]$ws.Range('CB15').Comment.Text('redox potential, measured relative to a hydrogen cell, indicating oxidation or reduction potential')
[void]$ws.Range('CC15').Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL')
[void]$ws.Range('CD15').Comment.Text('salinity measurement')
[void]$ws.Range('CE15').Comment.Text('Method or device employed for collecting sample')
[void]$ws.Range('CF15').Comment.Text('Processing applied to the sample during or after isolation')
[void]$ws.Range('CG15').Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
[void]$ws.Range('CH15').Comment.Text('duration for which sample was stored')
[void]$ws.Range('CI15').Comment.Text('location at which sample was stored, usually name of a specific freezer/room')
[void]$ws.Range('CJ15').Comment.Text('temperature at which sample was stored, e.g. -80')
[void]$ws.Range('CK15').Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
[void]$ws.Range('CL15').Comment.Text('information about the sediment type based on major constituents')
[void]$ws.Range('CM15').Comment.Text('concentration of silicate')
[void]$ws.Range('CN15').Comment.Text('sodium concentration')
[void]$ws.Range('CO15').Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
[void]$ws.Range('CP15').Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier')
[void]$ws.Range('CQ15').Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
[void]$ws.Range('CR15').Comment.Text('concentration of sulfate')
[void]$ws.Range('CS15').Comment.Text('concentration of sulfide')
[void]$ws.Range('CT15').Comment.Text('temperature of the sample at time of sampling')
[void]$ws.Range('CU15').Comment.Text('stage of tide')
[void]$ws.Range('CV15').Comment.Text('total carbon content')
[void]$ws.Range('CW15').Comment.Text('total nitrogen content of the sample')
[void]$ws.Range('CX15').Comment.Text('Definition for soil: total organic C content of the soil units of g C/kg soil. Definition otherwise: total organic carbon content')
[void]$ws.Range('CY15').Comment.Text('Feeding position in food chain (eg., chemolithotroph)')
[void]$ws.Range('CZ15').Comment.Text('turbidity measurement')
[void]$ws.Range('DA15').Comment.Text('water content measurement')

# Drop the now-duplicated trailing comment left over at DB15.
[void]$ws.Range('DB15').Comment.Delete()
